$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab18")

$ws.Range("O3").Value = 45.025733623187897
$ws.Range("O4").Value = 41.891955685351803
$ws.Range("O5").Value = 45.040888826404696
$ws.Range("O6").Value = 44.345664757404002
$ws.Range("O7").Value = 10.3727474306876
$ws.Range("O8").Value = 40.5940552908055
$ws.Range("O9").Value = 34.331492876482699
$ws.Range("O10").Value = 31.193814062404901
$ws.Range("O11").Value = 52.951540056276201
$ws.Range("O12").Value = 18.279911392911501
$ws.Range("O13").Value = 34.5944064075684
$ws.Range("O14").Value = 8.1227926531516701
$ws.Range("O15").Value = 16.4416262007732
$ws.Range("O16").Value = 13.688864979750299
$ws.Range("O17").Value = 34.529957898005001
$ws.Range("O18").Value = 72.891877155626503
$ws.Range("O19").Value = 39.553458604219401
$ws.Range("O20").Value = 44.622731725259897
$ws.Range("O21").Value = 39.862196110500797
$ws.Range("O22").Value = 15.309847129227901
$ws.Range("O23").Value = 33.056170556374198
$ws.Range("O24").Value = 11.516798569849399
$ws.Range("O25").Value = 149.939648857761
$ws.Range("O26").Value = 38.753396744285503
$ws.Range("O27").Value = 8.5738254450976097
$ws.Range("O28").Value = 10.7162886540966
$ws.Range("O29").Value = 23.032116454966602
$ws.Range("O30").Value = 27.8335084103309
$ws.Range("O31").Value = 19.064319772189801
$ws.Range("O32").Value = 120.178646783612
$ws.Range("O33").Value = 17.178050909260499
$ws.Range("O34").Value = 71.500734394283
$ws.Range("O35").Value = 17.0172856140736
$ws.Range("O36").Value = 14.1751950000597
$ws.Range("O37").Value = 14.385308686659201
$ws.Range("O38").Value = 14.6283437020997
$ws.Range("O39").Value = 25.605455733878198
$ws.Range("O40").Value = 10.5531691004687
$ws.Range("O42").Value = 33.2069791354443
$ws.Range("O43").Value = 32.985864623713397
$ws.Range("O44").Value = 42.439990494046803
$ws.Range("O45").Value = 18.783189127890399
$ws.Range("O46").Value = 23.482995586326499
$ws.Range("O47").Value = 29.4574859455408
$ws.Range("O48").Value = 28.9751184139811
$ws.Range("O49").Value = 22.632864580266698
$ws.Range("O50").Value = 9.1788319187919996
$ws.Range("O51").Value = 30.1948244523951
$ws.Range("O52").Value = 63.545103442808397
$ws.Range("O53").Value = 19.4231713892963
$ws.Range("O54").Value = 32.0125450422886
$ws.Range("O55").Value = 27.3971836887204
$ws.Range("O56").Value = 15.496881969881301
$ws.Range("O57").Value = 11.499661296247799
$ws.Range("O58").Value = 23.373588354331201
$ws.Range("O59").Value = 26.864112802491
$ws.Range("O60").Value = 22.5931424248264
$ws.Range("O61").Value = 17.668744907026699
$ws.Range("O62").Value = 21.428809165303299
$ws.Range("O63").Value = 29.074255902749801
$ws.Range("O64").Value = 27.112143836084599
$ws.Range("O65").Value = 24.0965283884358
$ws.Range("O66").Value = 28.688994384145701
$ws.Range("O67").Value = 15.827399720297199
$ws.Range("O68").Value = 17.5819635585536
$ws.Range("O69").Value = 13.418753152565801
$ws.Range("O70").Value = 36.685339214373897
$ws.Range("O71").Value = 17.668744907026699
$ws.Range("O72").Value = 13.3876320778232
$ws.Range("O73").Value = 32.443122353982901
$ws.Range("O74").Value = 30.495426343310001
$ws.Range("O75").Value = 43.522693584562703
$ws.Range("O76").Value = 52.997166297673402
$ws.Range("O77").Value = 21.190808658086699
$ws.Range("O78").Value = 50.405171963548
$ws.Range("O79").Value = 30.026167846821199
$ws.Range("O80").Value = 20.986150259167299
$ws.Range("O81").Value = 37.656844581745602
$ws.Range("O82").Value = 21.601179840296101
$ws.Range("O83").Value = 28.295846675813898
$ws.Range("O84").Value = 21.1779638950382
$ws.Range("O85").Value = 12.080169759101199
$ws.Range("O86").Value = 19.234074729172399
$ws.Range("O87").Value = 25.213955396743401
$ws.Range("O88").Value = 32.2537412319464
$ws.Range("O89").Value = 26.486034597540002
$ws.Range("O90").Value = 32.383493587139697
$ws.Range("O91").Value = 25.363513289198099
$ws.Range("O92").Value = 15.369241775483101
$ws.Range("O93").Value = 32.089359681950697
$ws.Range("O94").Value = 117.52915593533599
$ws.Range("O95").Value = 20.347630778781198
$ws.Range("O96").Value = 31.3023002992424
$ws.Range("O97").Value = 19.850907550630101
$ws.Range("O98").Value = 21.334067301334301
